# Powerpoint writer: consolidate text runs when possible.
# Slide 1 has two placeholders whose text was split across many single
# word/space <a:r> runs. Re-writing the text collapses the split runs
# into a single run. Because the flattened text content does not
# actually change, a plain TextRange.Text assignment is treated as a
# same-value write and skipped, so we go through a Characters()
# sub-range (whose assignment always rewrites the backing runs) instead.
# For the subtitle, the two line breaks (<a:br/>) in the middle of the
# paragraph must be preserved, so we update the text on either side of
# them via separate Characters() sub-ranges rather than the whole
# TextRange (which would turn embedded break characters into new
# paragraphs instead of <a:br/> soft breaks).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1: "Title 1" -> "Testing custom properties"
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Characters(1, $titleRange.Length).Text = "Testing custom properties"

# Shape 2: "Subtitle 2" -> "This is a subtitle" <br/><br/> "A. M."
$subRange = $s.Shapes.Item(2).TextFrame.TextRange

$firstPart = $subRange.Characters(1, 18)
$firstPart.Text = "This is a subtitle"

$secondPart = $subRange.Characters(21, 5)
$secondPart.Text = "A. M."
